# Fruta / hortaliza, semanal
# Update columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), P (Precio $/Kg) for rows 2-22 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44330
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 11).Value = 13000
$ws.Cells.Item(2, 12).Value = 14000
$ws.Cells.Item(2, 13).Value = 13500
$ws.Cells.Item(2, 16).Value = 1350
$ws.Cells.Item(3, 4).Value = 44263
$ws.Cells.Item(3, 10).Value = 300
$ws.Cells.Item(3, 11).Value = 15000
$ws.Cells.Item(3, 12).Value = 16000
$ws.Cells.Item(3, 13).Value = 15500
$ws.Cells.Item(3, 16).Value = 1550
$ws.Cells.Item(4, 4).Value = 44524
$ws.Cells.Item(4, 10).Value = 200
$ws.Cells.Item(4, 11).Value = 20000
$ws.Cells.Item(4, 12).Value = 21000
$ws.Cells.Item(4, 13).Value = 20500
$ws.Cells.Item(4, 16).Value = 2050
$ws.Cells.Item(5, 4).Value = 44694
$ws.Cells.Item(5, 10).Value = 400
$ws.Cells.Item(5, 11).Value = 16000
$ws.Cells.Item(5, 12).Value = 17000
$ws.Cells.Item(5, 13).Value = 16500
$ws.Cells.Item(5, 16).Value = 1650
$ws.Cells.Item(6, 4).Value = 44291
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(6, 11).Value = 13000
$ws.Cells.Item(6, 12).Value = 14000
$ws.Cells.Item(6, 13).Value = 13500
$ws.Cells.Item(6, 16).Value = 1350
$ws.Cells.Item(7, 4).Value = 44441
$ws.Cells.Item(7, 10).Value = 300
$ws.Cells.Item(7, 11).Value = 15000
$ws.Cells.Item(7, 12).Value = 16000
$ws.Cells.Item(7, 13).Value = 15500
$ws.Cells.Item(7, 16).Value = 1550
$ws.Cells.Item(8, 4).Value = 44377
$ws.Cells.Item(8, 10).Value = 650
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 15000
$ws.Cells.Item(8, 13).Value = 14538
$ws.Cells.Item(8, 16).Value = 1454
$ws.Cells.Item(9, 4).Value = 44265
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 15000
$ws.Cells.Item(9, 12).Value = 16000
$ws.Cells.Item(9, 13).Value = 15500
$ws.Cells.Item(9, 16).Value = 1550
$ws.Cells.Item(10, 4).Value = 44218
$ws.Cells.Item(10, 10).Value = 320
$ws.Cells.Item(10, 11).Value = 10000
$ws.Cells.Item(10, 12).Value = 11000
$ws.Cells.Item(10, 13).Value = 10500
$ws.Cells.Item(10, 16).Value = 1050
$ws.Cells.Item(11, 4).Value = 44160
$ws.Cells.Item(11, 10).Value = 360
$ws.Cells.Item(11, 11).Value = 10000
$ws.Cells.Item(11, 12).Value = 11000
$ws.Cells.Item(11, 13).Value = 10500
$ws.Cells.Item(11, 16).Value = 1050
$ws.Cells.Item(12, 4).Value = 44204
$ws.Cells.Item(12, 10).Value = 400
$ws.Cells.Item(12, 11).Value = 10000
$ws.Cells.Item(12, 12).Value = 11000
$ws.Cells.Item(12, 13).Value = 10500
$ws.Cells.Item(12, 16).Value = 1050
$ws.Cells.Item(13, 4).Value = 44679
$ws.Cells.Item(13, 10).Value = 200
$ws.Cells.Item(13, 11).Value = 19000
$ws.Cells.Item(13, 12).Value = 20000
$ws.Cells.Item(13, 13).Value = 19500
$ws.Cells.Item(13, 16).Value = 1950
$ws.Cells.Item(14, 4).Value = 44460
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 15000
$ws.Cells.Item(14, 12).Value = 16000
$ws.Cells.Item(14, 13).Value = 15500
$ws.Cells.Item(14, 16).Value = 1550
$ws.Cells.Item(15, 4).Value = 44580
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 18000
$ws.Cells.Item(15, 12).Value = 20000
$ws.Cells.Item(15, 13).Value = 19000
$ws.Cells.Item(15, 16).Value = 1900
$ws.Cells.Item(16, 4).Value = 44727
$ws.Cells.Item(16, 10).Value = 400
$ws.Cells.Item(16, 11).Value = 18000
$ws.Cells.Item(16, 12).Value = 19000
$ws.Cells.Item(16, 13).Value = 18500
$ws.Cells.Item(16, 16).Value = 1850
$ws.Cells.Item(17, 4).Value = 44644
$ws.Cells.Item(17, 10).Value = 300
$ws.Cells.Item(17, 11).Value = 20000
$ws.Cells.Item(17, 12).Value = 21000
$ws.Cells.Item(17, 13).Value = 20500
$ws.Cells.Item(17, 16).Value = 2050
$ws.Cells.Item(18, 4).Value = 44358
$ws.Cells.Item(18, 10).Value = 300
$ws.Cells.Item(18, 11).Value = 14000
$ws.Cells.Item(18, 12).Value = 15000
$ws.Cells.Item(18, 13).Value = 14500
$ws.Cells.Item(18, 16).Value = 1450
$ws.Cells.Item(19, 4).Value = 44428
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 15000
$ws.Cells.Item(19, 12).Value = 16000
$ws.Cells.Item(19, 13).Value = 15500
$ws.Cells.Item(19, 16).Value = 1550
$ws.Cells.Item(20, 4).Value = 44714
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 19000
$ws.Cells.Item(20, 12).Value = 20000
$ws.Cells.Item(20, 13).Value = 19500
$ws.Cells.Item(20, 16).Value = 1950
$ws.Cells.Item(21, 4).Value = 44406
$ws.Cells.Item(21, 10).Value = 400
$ws.Cells.Item(21, 11).Value = 14000
$ws.Cells.Item(21, 12).Value = 15000
$ws.Cells.Item(21, 13).Value = 14500
$ws.Cells.Item(21, 16).Value = 1450
$ws.Cells.Item(22, 4).Value = 44547
$ws.Cells.Item(22, 10).Value = 300
$ws.Cells.Item(22, 11).Value = 19000
$ws.Cells.Item(22, 12).Value = 20000
$ws.Cells.Item(22, 13).Value = 19500
$ws.Cells.Item(22, 16).Value = 1950